# Build site at 2023-04-12 14:53:07 UTC
# Update LOQ4094.xlsx: fill in the Petrochemical course syllabus content
# (objectives, summary/full programs, method/criteria text, bibliography)
# and insert the missing "Docentes responsaveis" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- long text blocks reused below -----------------------------------
$objetivos = 'Introdução à indústria petroquímica, a partir da obtenção de matérias-primas básicas até a fabricação de produtos de segunda geração.'

$programaResumido = '1.Fundamentos da indústria petroquímica: interface refino-petroquímica, matérias-primas, cadeia industrial petroquímica, complexos petroquímicos, petroquímica brasileira;2.Produção e separação de olefinas: craqueamento a vapor, licenciadores de tecnologia, unidades de cracking, corte C4;3.Produção e Separação de Aromáticos: matéria-prima, reforma catalítica, licenciadores, processo de extração dos aromáticos, separação de BTX;4.Produção e utilização do gás de síntese: processos de produção, amônia, metanol;5.Produção de intermediários petroquímicos: etilbenzeno, estireno, cumeno, fenol, ácido tereftálico; óxido de eteno, intermediários para fibras sintéticas;6.Produção de polímeros sintéticos: polietileno, polipropileno, polímeros vinílicos, poliestireno, poli (tereftalato de etileno); poliamidas, poliuretanos, elastômeros'

$programaCompleto = 'Fundamentos da indústria petroquímica: interface refino-petroquímica, matérias-primas, cadeia industrial petroquímica, complexos petroquímicos, petroquímica brasileira; Produção e separação de olefinas: craqueamento a vapor, licenciadores de tecnologia, unidades de cracking, corte C4; Produção e Separação de Aromáticos: matéria-prima, reforma catalítica, licenciadores, processo de extração dos aromáticos, separação de BTX; Produção e utilização do gás de síntese: processos de produção, amônia, metanol; produção de intermediários petroquímicos: etilbenzeno, estireno, cumeno, fenol, ácido tereftálico; óxido de eteno, intermediários para fibras sintéticas; Produção de polímeros sintéticos: polietileno, polipropileno, polímeros vinílicos, poliestireno, poli(tereftalato de eteno; poliamidas, poliuretanos, elastômeros.'

$metodo = 'Aulas expositivas, desenvolvimento de exercícios em sala e fora de sala de aula, discussão de casos práticos.'

$criterio = 'Provas, avaliação através de exercícios ou casos práticos elaborados fora de sala de aula.'

$normaRecuperacao = 'Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação.'

$bibliografia = 'a)Meyers, R. A., Handbook of Petrochemicals Production Process, The McGraw Hill Companies, 1ª Edição, 2005;
b)Speight, J. G., The Chemistry and Technology of Petroleum, CRC Press, 4ª Edição, 2007;
c)Perrone, O. V., Silva Filho, A. P. (Coordenadores), Processos Petroquímicos, Editora Sinergia, 1ª Edição, 2013;
d)Leite, L. F., Olefinas Leves, Editora Interciência, 1ª edição, 2012.
e)Brasil, N. I., Araújo, M. A. S., Souza, E. C. M, Processamento de Petróleo e Gás, Editora LTC, 1ª Edição, 2012;
f)Fundamentos do Refino do Petróleo  Tecnologia e Economia, Szklo, A. S., Uller, V. C., Bonfá, M. H. P., Editora Interciência, 3ª Edição, 2012;
g)Oil and Gas Journal;
h)Revista Petro & Química.'

$professor = '1285870 - Marcos Villela Barcza'

# --- insert the missing "Docentes responsaveis:" data row ------------
# Row 13 ("Programa resumido:") shifts down to make room for the
# professor name that belongs next to "Docentes responsaveis:" (row 12).
$ws.Rows.Item(13).Insert()

# The blank row Excel creates inherits column-A formatting in A13 and
# copies column A's style into B13 - clear A13 entirely and repair B13/C13
# formatting by pulling the (correct) formats down from row 14 before
# writing the real values.
$ws.Range("A13").Clear()
$ws.Range("B14:C14").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Objetivos: (row 10) now gets the real course description --------
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# --- Docentes responsaveis: (row 13, newly inserted) ------------------
$ws.Range("B13").Value = $professor
$ws.Range("C13").Value = $professor

# --- Programa resumido: (row 14) --------------------------------------
$ws.Range("B14").Value = $programaResumido
$ws.Range("C14").Value = $programaResumido

# --- Short syllabus: (row 15) gains the 60pt custom row height --------
$ws.Rows.Item(15).RowHeight = 60

# --- Programa: (row 16) ------------------------------------------------
$ws.Range("B16").Value = $programaCompleto
$ws.Range("C16").Value = $programaCompleto

# --- Metodo: (row 19) ---------------------------------------------------
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# --- Criterio: (row 20) --------------------------------------------------
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# --- Norma de recuperacao: (row 21) ---------------------------------------
$ws.Range("B21").Value = $normaRecuperacao
$ws.Range("C21").Value = $normaRecuperacao

# --- Bibliografia: (row 22) -----------------------------------------------
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia
